$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sfrp2"
$ws.Range("C2").Value = "Fzd5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.3709053333333334
$ws.Range("H2").Value = 1.112716
$ws.Range("I2").Value = 0.01868139080022812
$ws.Range("J2").Value = 0.01868139080022812
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 10.000565
$ws.Range("N2").Value = 30.001695
$ws.Range("O2").Value = 0.6316353758144477
$ws.Range("P2").Value = 0.6316353758144477
$ws.Range("Q2").Value = 3.709262894846667
$ws.Range("R2").Value = 33.38336605362
$ws.Range("S2").Value = 0.01179982729883866
$ws.Range("T2").Value = 0.01179982729883866

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sfrp2"
$ws.Range("C3").Value = "Fzd5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.3709053333333334
$ws.Range("H3").Value = 1.112716
$ws.Range("I3").Value = 0.01868139080022812
$ws.Range("J3").Value = 0.01868139080022812
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.264793333333333
$ws.Range("N3").Value = 12.79438
$ws.Range("O3").Value = 0.2693642149089528
$ws.Range("P3").Value = 0.2693642149089528
$ws.Range("Q3").Value = 1.581834592897778
$ws.Range("R3").Value = 14.23651133608
$ws.Range("S3").Value = 0.005032098166310782
$ws.Range("T3").Value = 0.005032098166310782

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sfrp2"
$ws.Range("C4").Value = "Fzd5"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.3709053333333334
$ws.Range("H4").Value = 1.112716
$ws.Range("I4").Value = 0.01868139080022812
$ws.Range("J4").Value = 0.01868139080022812
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.567455
$ws.Range("N4").Value = 4.702364999999999
$ws.Range("O4").Value = 0.09900040927659938
$ws.Range("P4").Value = 0.09900040927659938
$ws.Range("Q4").Value = 0.58137741926
$ws.Range("R4").Value = 5.23239677334
$ws.Range("S4").Value = 0.001849465335078682
$ws.Range("T4").Value = 0.001849465335078682

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Sfrp2"
$ws.Range("C5").Value = "Fzd5"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 18.59812133333334
$ws.Range("H5").Value = 55.794364
$ws.Range("I5").Value = 0.936731671274772
$ws.Range("J5").Value = 0.936731671274772
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10.000565
$ws.Range("N5").Value = 30.001695
$ws.Range("O5").Value = 0.6316353758144477
$ws.Range("P5").Value = 0.6316353758144477
$ws.Range("Q5").Value = 185.9917212718867
$ws.Range("R5").Value = 1673.92549144698
$ws.Range("S5").Value = 0.5916728612229363
$ws.Range("T5").Value = 0.5916728612229363

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Sfrp2"
$ws.Range("C6").Value = "Fzd5"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 18.59812133333334
$ws.Range("H6").Value = 55.794364
$ws.Range("I6").Value = 0.936731671274772
$ws.Range("J6").Value = 0.936731671274772
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.264793333333333
$ws.Range("N6").Value = 12.79438
$ws.Range("O6").Value = 0.2693642149089528
$ws.Range("P6").Value = 0.2693642149089528
$ws.Range("Q6").Value = 79.31714387492445
$ws.Range("R6").Value = 713.8542948743201
$ws.Range("S6").Value = 0.2523219912132802
$ws.Range("T6").Value = 0.2523219912132802

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Sfrp2"
$ws.Range("C7").Value = "Fzd5"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 18.59812133333334
$ws.Range("H7").Value = 55.794364
$ws.Range("I7").Value = 0.936731671274772
$ws.Range("J7").Value = 0.936731671274772
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.567455
$ws.Range("N7").Value = 4.702364999999999
$ws.Range("O7").Value = 0.09900040927659938
$ws.Range("P7").Value = 0.09900040927659938
$ws.Range("Q7").Value = 29.15171827454
$ws.Range("R7").Value = 262.36546447086
$ws.Range("S7").Value = 0.09273681883855538
$ws.Range("T7").Value = 0.09273681883855538

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sfrp2"
$ws.Range("C8").Value = "Fzd5"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8852410000000001
$ws.Range("H8").Value = 2.655723
$ws.Range("I8").Value = 0.04458693792499994
$ws.Range("J8").Value = 0.04458693792499994
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.000565
$ws.Range("N8").Value = 30.001695
$ws.Range("O8").Value = 0.6316353758144477
$ws.Range("P8").Value = 0.6316353758144477
$ws.Range("Q8").Value = 8.852910161165001
$ws.Range("R8").Value = 79.67619145048499
$ws.Range("S8").Value = 0.02816268729267279
$ws.Range("T8").Value = 0.02816268729267279

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sfrp2"
$ws.Range("C9").Value = "Fzd5"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8852410000000001
$ws.Range("H9").Value = 2.655723
$ws.Range("I9").Value = 0.04458693792499994
$ws.Range("J9").Value = 0.04458693792499994
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 4.264793333333333
$ws.Range("N9").Value = 12.79438
$ws.Range("O9").Value = 0.2693642149089528
$ws.Range("P9").Value = 0.2693642149089528
$ws.Range("Q9").Value = 3.775369915193334
$ws.Range("R9").Value = 33.97832923674
$ws.Range("S9").Value = 0.01201012552936182
$ws.Range("T9").Value = 0.01201012552936182

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sfrp2"
$ws.Range("C10").Value = "Fzd5"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8852410000000001
$ws.Range("H10").Value = 2.655723
$ws.Range("I10").Value = 0.04458693792499994
$ws.Range("J10").Value = 0.04458693792499994
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.567455
$ws.Range("N10").Value = 4.702364999999999
$ws.Range("O10").Value = 0.09900040927659938
$ws.Range("P10").Value = 0.09900040927659938
$ws.Range("Q10").Value = 1.387575431655
$ws.Range("R10").Value = 12.488178884895
$ws.Range("S10").Value = 0.004414125102965324
$ws.Range("T10").Value = 0.004414125102965324
